# [Kadastro App] Yeni kayit eklendi: 2969
#
# Appends one new record row (Kayit No 2969, Erdemli birimi, 2025-09-10
# tarihli CAP isi, AYHAN KARADAYI teknisyeni) to the bottom of both the
# "Kayitlar" master log sheet and the "Erdemli" birim sheet, which mirror
# the same records. All six columns in this workbook are stored as text
# (even the numeric-looking "Kayit No"/"Parsel Sayisi" and the date), so
# we force a text number format before writing the numeric/date-looking
# values to stop them being auto-coerced into a number/date, then restore
# the cell style so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook

$newRecord = @{
    1 = "2969"
    2 = "2025-09-10"
    3 = "Erdemli"
    4 = "1"
    5 = "ÇAP"
    6 = "AYHAN KARADAYI (K.Teknisyeni)"
}

# Columns whose values look numeric/date-like and therefore need to be
# pinned to text so Excel doesn't silently convert them to a number/date
# (mirrors the existing rows, which are all stored as literal text).
$textLikeColumns = @(1, 2, 4)

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = $ws.UsedRange.Rows.Count + 1

    foreach ($col in 1..6) {
        $cell = $ws.Cells.Item($newRow, $col)

        if ($textLikeColumns -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $newRecord[$col]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newRecord[$col]
        }
    }

    # Keep the "numbers stored as text" ignored-error range in sync with
    # the newly extended data range (best effort - some hosts don't
    # surface this as a settable property).
    try {
        $ws.Range("A1:F$newRow").Errors.Item(9).Ignore = $true
    } catch {
    }
}
